$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values (rows 2-11)
$ws.Range("C2").Value = 2670005915.8315001
$ws.Range("C3").Value = 2339240926
$ws.Range("C4").Value = 1138443110
$ws.Range("C5").Value = 237165911
$ws.Range("C6").Value = 211681494
$ws.Range("C7").Value = 100111524
$ws.Range("C8").Value = 63548766
$ws.Range("C9").Value = 56996640
$ws.Range("C10").Value = 45142346
$ws.Range("C11").Value = 51573653

# Update selection to C4
$ws.Range("C4").Select()
